# Literature overview - Systems of Systems additions + review feedback
# - Widen column C (Scope) to fit "Theoretical Background"
# - Fill in two new literature rows (AMS053 / AMS054) that previously only had an ID
# - Reposition the view/selection to the newly-edited area

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlTop = -4160

# --- Column C width (Scope) ---------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 25.3

# --- Row 54 : AMS053 - Maier, M. - Architecting Principles for Systems-of-Systems
$ws.Range("B54").Value = "Systems of Systems Definition"
$ws.Range("B54").VerticalAlignment = $xlTop

$ws.Range("C54").Value = "Theoretical Background"
$ws.Range("C54").VerticalAlignment = $xlTop

$ws.Range("D54").Value = "Architecting Principles for Systems-of-Systems"
$ws.Range("D54").VerticalAlignment = $xlTop
$ws.Range("D54").WrapText = $true

$ws.Range("F54").Value = "Maier, M."
$ws.Range("F54").VerticalAlignment = $xlTop
$ws.Range("F54").WrapText = $true

$ws.Range("G54").Value = 1996
$ws.Range("G54").VerticalAlignment = $xlTop

$ws.Range("H54").Value = "Journal Article"
$ws.Range("H54").VerticalAlignment = $xlTop

$ws.Range("I54").Value = "Maier1996"
$ws.Range("I54").VerticalAlignment = $xlTop

$ws.Range("J54").Value = "Yes"
$ws.Range("J54").VerticalAlignment = $xlTop

$ws.Range("K54").Value = "Yes"
$ws.Range("K54").VerticalAlignment = $xlTop

$ws.Range("L54").Value = "Yes"
$ws.Range("L54").VerticalAlignment = $xlTop

$ws.Range("M54").Value = "Semantic Scholar"
$ws.Range("M54").VerticalAlignment = $xlTop

$ws.Range("N54").Value = "10.1002/J.2334-5837.1996.TB02054.X"
$ws.Range("N54").VerticalAlignment = $xlTop

$ws.Range("O54").Value = "https://doi.org/10.1002/J.2334-5837.1996.TB02054.X"
$ws.Range("O54").VerticalAlignment = $xlTop
$ws.Range("O54").WrapText = $true

$ws.Range("P54").Value = "September,2021"
$ws.Range("P54").VerticalAlignment = $xlTop

$ws.Range("Q54").Value = "Read"
$ws.Range("Q54").VerticalAlignment = $xlTop

$ws.Range("R54").Value = "No"
$ws.Range("R54").VerticalAlignment = $xlTop

$ws.Range("T54").Value = "Defining Systems of Systems"
$ws.Range("T54").VerticalAlignment = $xlTop

$ws.Range("U54").Value = "Old paper but a lot of other sources refer to the “Maier’s criteri"
$ws.Range("U54").VerticalAlignment = $xlTop

$ws.Rows.Item(54).RowHeight = 30

# --- Row 55 : AMS054 - Dersin, Pierre - IEEE Systems of Systems Whitepaper
$ws.Range("B55").Value = "Systems of Systems Definition"
$ws.Range("B55").VerticalAlignment = $xlTop

$ws.Range("C55").Value = "Theoretical Background"
$ws.Range("C55").VerticalAlignment = $xlTop

$ws.Range("D55").Value = "IEEE Systems of Systems Whitepaper"
$ws.Range("D55").VerticalAlignment = $xlTop
$ws.Range("D55").WrapText = $true

$ws.Range("F55").Value = "Dersin, Pierre"
$ws.Range("F55").VerticalAlignment = $xlTop
$ws.Range("F55").WrapText = $true

$ws.Range("G55").Value = 2014
$ws.Range("G55").VerticalAlignment = $xlTop

$ws.Range("H55").Value = "Whitepaper"
$ws.Range("H55").VerticalAlignment = $xlTop

$ws.Range("I55").Value = "Dersin2014"
$ws.Range("I55").VerticalAlignment = $xlTop

$ws.Range("J55").Value = "Yes"
$ws.Range("J55").VerticalAlignment = $xlTop

$ws.Range("K55").Value = "N/A"
$ws.Range("K55").VerticalAlignment = $xlTop

$ws.Range("L55").Value = "Yes"
$ws.Range("L55").VerticalAlignment = $xlTop

$ws.Range("M55").Value = "IEEE"
$ws.Range("M55").VerticalAlignment = $xlTop

$ws.Range("O55").Value = "https://rs.ieee.org/technical-activities/technical-committees/systems-of-systems.html"
$ws.Range("O55").VerticalAlignment = $xlTop
$ws.Range("O55").WrapText = $true

$ws.Range("P55").Value = "September,2021"
$ws.Range("P55").VerticalAlignment = $xlTop

$ws.Range("Q55").Value = "Read"
$ws.Range("Q55").VerticalAlignment = $xlTop

$ws.Range("R55").Value = "No"
$ws.Range("R55").VerticalAlignment = $xlTop

$ws.Range("T55").Value = "Defining Systems of Systems as an extension on Maier (AMS053)"
$ws.Range("T55").VerticalAlignment = $xlTop
$ws.Range("T55").WrapText = $true

$ws.Rows.Item(55).RowHeight = 60

# --- View state: land on the newly-edited row -----------------------------
$ws.Range("A54").Select()
